$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated odds for row 5
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 4.45
$ws.Range("J5").Value = 6.7
$ws.Range("K5").Value = 2.37
$ws.Range("L5").Value = 1.85
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 8.25
$ws.Range("O5").Value = 1.24
$ws.Range("P5").Value = 3.65
$ws.Range("Q5").Value = 1.72
$ws.Range("R5").Value = 2.05
$ws.Range("S5").Value = 1.35
$ws.Range("T5").Value = 2.95
$ws.Range("W5").Value = 18
$ws.Range("X5").Value = 45
$ws.Range("Y5").Value = 23
$ws.Range("Z5").Value = 175
$ws.Range("AC5").Value = 8.25
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 22
$ws.Range("AG5").Value = 6.7
$ws.Range("AH5").Value = 6.3
$ws.Range("AN5").Value = 8.25
$ws.Range("AQ5").Value = 300
$ws.Range("AT5").Value = 2.95
$ws.Range("AX5").Value = 6.2
$ws.Range("AY5").Value = 17
$ws.Range("AZ5").Value = 17.5

# Updated odds for row 7
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 1.82
$ws.Range("J7").Value = 4.5
$ws.Range("K7").Value = 2.12
$ws.Range("L7").Value = 2.4
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 6.8
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.75
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.72
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.85
$ws.Range("W7").Value = 10.75
$ws.Range("X7").Value = 22
$ws.Range("Y7").Value = 13.5
$ws.Range("Z7").Value = 65
$ws.Range("AA7").Value = 40
$ws.Range("AB7").Value = 45
$ws.Range("AC7").Value = 6.8
$ws.Range("AD7").Value = 6.5
$ws.Range("AE7").Value = 15.5
$ws.Range("AF7").Value = 90
$ws.Range("AG7").Value = 6.6
$ws.Range("AH7").Value = 8.25
$ws.Range("AJ7").Value = 15
$ws.Range("AK7").Value = 15
$ws.Range("AL7").Value = 29
$ws.Range("AN7").Value = 5.9
$ws.Range("AO7").Value = 23
$ws.Range("AP7").Value = 29
$ws.Range("AQ7").Value = 120
$ws.Range("AT7").Value = 2.72
$ws.Range("AU7").Value = 7.3
$ws.Range("AV7").Value = 65
$ws.Range("AW7").Value = 3.7
$ws.Range("AX7").Value = 9.25
$ws.Range("AZ7").Value = 32
$ws.Range("BA7").Value = 65
